$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = -213109876.32717904
    "C2" = -51642405.72691207
    "D2" = 109825064.87335604
    "E2" = 271292535.4736225
    "F2" = 432760006.07389426
    "B3" = -103282078.75656138
    "C3" = 58185391.84370559
    "D3" = 219652862.4439737
    "E3" = 381120333.0442401
    "F3" = 542587803.6445119
    "B4" = 116381714.73266701
    "C4" = 277849185.332934
    "D4" = 439316655.9332021
    "E4" = 600784126.5334686
    "F4" = 762251597.1337404
    "B5" = 379992696.6208639
    "C5" = 541460167.221131
    "D5" = 702927637.821399
    "E5" = 864395108.4216654
    "F5" = 1025862579.0219373
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
